$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "noSQL options" entry text
$ws.Range("E9").Value = "Looking into noSQL options for this app and fending off error after error trying to implement them. Implemented basis of Hive"

# Update the project repo URL
$ws.Range("E1").Value = "https://github.com/leono93/world-of-warcraft-pocket-buddy"

# Fill in previously empty row 10 (start/end times + description)
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 14
$ws.Range("E10").Value = "Reading Blizzard API documentation and public resources on it as well as NodeJS implementation and libraries"

# Fill in previously empty row 11 (start/end times)
$ws.Range("B11").Value = 11
$ws.Range("C11").Value = 16

# Move the active selection to G13
$ws.Range("G13").Select()
